$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B4:G8").ClearFormats()
